$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 144.42857
$ws.Range("I6").Value = 132.2
$ws.Range("J6").Value = 175
$ws.Range("K6").Value = 396.6
$ws.Range("L6").Value = 525
$ws.Range("M6").Value = -284.6
$ws.Range("N6").Value = -749
$ws.Range("H28").Value = 6005
$ws.Range("I28").Value = 6005
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 6005
$ws.Range("L28").Value = 0
$ws.Range("M28").ClearContents()
$ws.Range("N28").Value = -5520
$ws.Range("H33").Value = 332.22223
$ws.Range("I33").Value = 242.46666
$ws.Range("K33").Value = 242.46666
$ws.Range("M33").Value = -13.46665999999999
$ws.Range("H64").Value = 9861.111000000001
$ws.Range("I64").Value = 10350
$ws.Range("J64").Value = 9250
$ws.Range("K64").Value = 10350
$ws.Range("L64").Value = 9250
$ws.Range("M64").Value = -10102
$ws.Range("N64").Value = -9746
$ws.Range("H67").Value = 9861.111000000001
$ws.Range("I67").Value = 10350
$ws.Range("J67").Value = 9250
$ws.Range("K67").Value = 10350
$ws.Range("L67").Value = 9250
$ws.Range("M67").Value = -9492
$ws.Range("N67").Value = -10966
$ws.Range("H70").Value = 23124
$ws.Range("J70").Value = 23124
$ws.Range("L70").Value = 69372
$ws.Range("N70").Value = -69912
$ws.Range("H73").Value = 23124
$ws.Range("J73").Value = 23124
$ws.Range("L73").Value = 69372
$ws.Range("N73").Value = -71244
$ws.Range("H74").Value = 5716.8335
$ws.Range("I74").Value = 4859.4
$ws.Range("K74").Value = 4859.4
$ws.Range("M74").Value = -3923.4
$ws.Range("H77").Value = 5716.8335
$ws.Range("I77").Value = 4859.4
$ws.Range("K77").Value = 24297
$ws.Range("M77").Value = -19617
$ws.Range("H80").Value = 826.2
$ws.Range("I80").Value = 635.25
$ws.Range("J80").Value = 1044.4286
$ws.Range("K80").Value = 1905.75
$ws.Range("L80").Value = 3133.2858
$ws.Range("M80").Value = -907.75
$ws.Range("N80").Value = -5129.2858
$ws.Range("H83").Value = 826.2
$ws.Range("I83").Value = 635.25
$ws.Range("J83").Value = 1044.4286
$ws.Range("K83").Value = 5717.25
$ws.Range("L83").Value = 9399.857399999999
$ws.Range("M83").Value = -725.25
$ws.Range("N83").Value = -19383.8574
$ws.Range("H112").Value = 2722.805
$ws.Range("I112").Value = 1422
$ws.Range("J112").Value = 2755.325
$ws.Range("K112").Value = 4266
$ws.Range("L112").Value = 8265.974999999999
$ws.Range("M112").Value = -3158
$ws.Range("N112").Value = -10481.975
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").ClearContents()
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = 0
$ws.Range("H127").Value = 999
$ws.Range("I127").Value = 999
$ws.Range("K127").Value = 2997
$ws.Range("M127").Value = 1963
$ws.Range("H135").Value = 3543.3157
$ws.Range("J135").Value = 2000
$ws.Range("L135").Value = 18000
$ws.Range("N135").Value = -23070
$ws.Range("H137").Value = 1760.5769
$ws.Range("I137").Value = 1632.2778
$ws.Range("J137").Value = 2049.25
$ws.Range("K137").Value = 4896.8334
$ws.Range("L137").Value = 6147.75
$ws.Range("M137").Value = -2346.8334
$ws.Range("N137").Value = -11247.75
$ws.Range("H138").Value = 3557.6316
$ws.Range("J138").Value = 4234.4546
$ws.Range("L138").Value = 12703.3638
$ws.Range("N138").Value = -22983.3638

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2016
$ws.Range("I61").Value = 1945.6111
$ws.Range("J61").Value = 2649.5
$ws.Range("K61").Value = 1945.6111
$ws.Range("L61").Value = 2649.5
$ws.Range("M61").Value = -1733.6111
$ws.Range("N61").Value = -3073.5
$ws.Range("H110").Value = 984.1053000000001
$ws.Range("I110").Value = 835.17645
$ws.Range("K110").Value = 835.17645
$ws.Range("M110").Value = 1209.82355
$ws.Range("H122").Value = 4768.9766
$ws.Range("I122").Value = 4596.2705
$ws.Range("K122").Value = 13788.8115
$ws.Range("M122").Value = -11338.8115
$ws.Range("H136").Value = 2016
$ws.Range("I136").Value = 1945.6111
$ws.Range("J136").Value = 2649.5
$ws.Range("K136").Value = 5836.8333
$ws.Range("L136").Value = 7948.5
$ws.Range("M136").Value = -3286.8333
$ws.Range("N136").Value = -13048.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1057.9445
$ws.Range("I22").Value = 416.25
$ws.Range("J22").Value = 2341.3333
$ws.Range("K22").Value = 416.25
$ws.Range("L22").Value = 2341.3333
$ws.Range("M22").Value = -66.25
$ws.Range("N22").Value = -3041.3333
$ws.Range("H31").Value = 1133.3636
$ws.Range("J31").Value = 998.5
$ws.Range("L31").Value = 998.5
$ws.Range("N31").Value = -1588.5
$ws.Range("H34").Value = 1133.3636
$ws.Range("J34").Value = 998.5
$ws.Range("L34").Value = 998.5
$ws.Range("N34").Value = -1402.5
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").ClearContents()
$ws.Range("N133").Value = 0

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 16668139
$ws.Range("J2").Value = 37039884
$ws.Range("L2").Value = 222239304
$ws.Range("N2").Value = -222239530
$ws.Range("H74").Value = 20000
$ws.Range("J74").Value = 20000
$ws.Range("L74").Value = 60000
$ws.Range("N74").Value = -62122
$ws.Range("H77").Value = 20000
$ws.Range("J77").Value = 20000
$ws.Range("L77").Value = 180000
$ws.Range("N77").Value = -190608
$ws.Range("H107").Value = 835.5
$ws.Range("I107").Value = 797
$ws.Range("K107").Value = 2391
$ws.Range("M107").Value = -471
$ws.Range("H122").Value = 9998.666999999999
$ws.Range("J122").Value = 9998.5
$ws.Range("L122").Value = 89986.5
$ws.Range("N122").Value = -94886.5
$ws.Range("H140").Value = 1493.5
$ws.Range("I140").Value = 1072.375
$ws.Range("J140").Value = 2055
$ws.Range("K140").Value = 3217.125
$ws.Range("L140").Value = 6165
$ws.Range("M140").Value = 1962.875
$ws.Range("N140").Value = -16525

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1452.1111
$ws.Range("I122").Value = 1483.5714
$ws.Range("J122").Value = 1342
$ws.Range("K122").Value = 4450.7142
$ws.Range("L122").Value = 4026
$ws.Range("M122").Value = -2000.7142
$ws.Range("N122").Value = -8926

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6052
$ws.Range("I7").Value = 4104.75
$ws.Range("J7").Value = 7999.25
$ws.Range("K7").Value = 4104.75
$ws.Range("L7").Value = 7999.25
$ws.Range("M7").Value = -3992.75
$ws.Range("N7").Value = -8223.25
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("H82").Value = 2393.2593
$ws.Range("I82").Value = 1774
$ws.Range("J82").Value = 3060.1538
$ws.Range("K82").Value = 1774
$ws.Range("L82").Value = 3060.1538
$ws.Range("M82").Value = -1413
$ws.Range("N82").Value = -3782.1538
$ws.Range("H85").Value = 2393.2593
$ws.Range("I85").Value = 1774
$ws.Range("J85").Value = 3060.1538
$ws.Range("K85").Value = 1774
$ws.Range("L85").Value = 3060.1538
$ws.Range("M85").Value = -526
$ws.Range("N85").Value = -5556.1538
$ws.Range("H93").Value = 2347.875
$ws.Range("I93").Value = 1493
$ws.Range("J93").Value = 4912.5
$ws.Range("K93").Value = 1493
$ws.Range("L93").Value = 4912.5
$ws.Range("M93").Value = -245
$ws.Range("N93").Value = -7408.5
$ws.Range("H103").Value = 12000
$ws.Range("J103").Value = 12000
$ws.Range("L103").Value = 12000
$ws.Range("N103").Value = -14344
$ws.Range("H126").Value = 6052
$ws.Range("I126").Value = 4104.75
$ws.Range("J126").Value = 7999.25
$ws.Range("K126").Value = 12314.25
$ws.Range("L126").Value = 23997.75
$ws.Range("M126").Value = -9844.25
$ws.Range("N126").Value = -28937.75
$ws.Range("H136").Value = 3747.4167
$ws.Range("I136").Value = 3596.6
$ws.Range("J136").Value = 4501.5
$ws.Range("K136").Value = 10789.8
$ws.Range("L136").Value = 13504.5
$ws.Range("M136").Value = -8239.799999999999
$ws.Range("N136").Value = -18604.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1712
$ws.Range("J81").Value = 1859.6666
$ws.Range("L81").Value = 3719.3332
$ws.Range("N81").Value = -5841.3332
$ws.Range("H84").Value = 1712
$ws.Range("J84").Value = 1859.6666
$ws.Range("L84").Value = 18596.666
$ws.Range("N84").Value = -29204.666
$ws.Range("H126").Value = 1325.3077
$ws.Range("I126").Value = 1043
$ws.Range("J126").Value = 2266.3333
$ws.Range("K126").Value = 3129
$ws.Range("L126").Value = 6798.999899999999
$ws.Range("M126").Value = -659
$ws.Range("N126").Value = -11738.9999

